$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-06-08 Sunday" "2025-06-09 Monday"

Replace-Text "81×80=" "56×42="
Replace-Text "64×35=" "68×18="
Replace-Text "97×51=" "23×42="
Replace-Text "15×44=" "65×55="
Replace-Text "28×26=" "33×98="

Replace-Text "68×32=" "32×19="
Replace-Text "33×66=" "81×24="
Replace-Text "54×64=" "57×63="
Replace-Text "94×39=" "84×74="
Replace-Text "79×38=" "41×72="

Replace-Text "43×88=" "33×20="
Replace-Text "75×99=" "20×67="
Replace-Text "90×21=" "16×98="
Replace-Text "39×60=" "26×13="
Replace-Text "36×24=" "47×64="

Replace-Text "28×91=" "45×16="
Replace-Text "91×96=" "77×24="
Replace-Text "59×55=" "88×33="
Replace-Text "12×37=" "32×58="
Replace-Text "90×34=" "55×61="

Replace-Text "93×81=" "93×87="
Replace-Text "92×99=" "87×69="
Replace-Text "23×66=" "25×22="
Replace-Text "59×79=" "53×26="
Replace-Text "72×73=" "14×22="
